# Update market-price-derived profit figures (H:N) across all job sheets
# per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2606.7856
$ws.Range("I100").Value = 2656.4285
$ws.Range("J100").Value = 2557.1428
$ws.Range("K100").Value = 2656.4285
$ws.Range("L100").Value = 2557.1428
$ws.Range("M100").Value = -2115.4285
$ws.Range("N100").Value = -3639.1428

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3956.4075
$ws.Range("I137").Value = 3707.7
$ws.Range("K137").Value = 11123.1
$ws.Range("M137").Value = -8573.099999999999

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 207291.78
$ws.Range("J138").Value = 287954.6
$ws.Range("L138").Value = 863863.7999999999
$ws.Range("N138").Value = -874143.7999999999

# ARM row 7
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 41000
$ws.Range("J7").Value = 41000
$ws.Range("L7").Value = 41000
$ws.Range("N7").Value = -41228

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 479548.97
$ws.Range("I32").Value = 642268.4
$ws.Range("J32").Value = 10534.177
$ws.Range("K32").Value = 642268.4
$ws.Range("L32").Value = 10534.177
$ws.Range("M32").Value = -641981.4
$ws.Range("N32").Value = -11108.177

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6346.364
$ws.Range("I132").Value = 8762.4
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 26287.2
$ws.Range("L132").Value = 12999
$ws.Range("M132").Value = -23757.2
$ws.Range("N132").Value = -18059

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2472.5
$ws.Range("I99").Value = 2472.5
$ws.Range("K99").Value = 2472.5
$ws.Range("M99").Value = -974.5

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3076.5
$ws.Range("I105").Value = 3013.2222
$ws.Range("J105").Value = 3266.3333
$ws.Range("K105").Value = 3013.2222
$ws.Range("L105").Value = 3266.3333
$ws.Range("M105").Value = -1266.2222
$ws.Range("N105").Value = -6760.3333

# BSM row 126
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 70000
$ws.Range("J126").Value = 70000
$ws.Range("L126").Value = 70000
$ws.Range("N126").Value = -79880

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 633.875
$ws.Range("I16").Value = 591.7778
$ws.Range("J16").Value = 688
$ws.Range("K16").Value = 591.7778
$ws.Range("L16").Value = 688
$ws.Range("M16").Value = -304.7778
$ws.Range("N16").Value = -1262

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2778.963
$ws.Range("I31").Value = 1019
$ws.Range("J31").Value = 5338.909
$ws.Range("K31").Value = 1019
$ws.Range("L31").Value = 5338.909
$ws.Range("M31").Value = -724
$ws.Range("N31").Value = -5928.909

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2778.963
$ws.Range("I34").Value = 1019
$ws.Range("J34").Value = 5338.909
$ws.Range("K34").Value = 1019
$ws.Range("L34").Value = 5338.909
$ws.Range("M34").Value = -817
$ws.Range("N34").Value = -5742.909

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 633.875
$ws.Range("I113").Value = 591.7778
$ws.Range("J113").Value = 688
$ws.Range("K113").Value = 591.7778
$ws.Range("L113").Value = 688
$ws.Range("M113").Value = 1578.2222
$ws.Range("N113").Value = -5028

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 18521102
$ws.Range("I132").Value = 1578.3334
$ws.Range("K132").Value = 4735.0002
$ws.Range("M132").Value = -2205.0002

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1166
$ws.Range("J5").Value = 1556.3334
$ws.Range("L5").Value = 4669.0002
$ws.Range("N5").Value = -4893.0002

# CUL row 76
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 500
$ws.Range("I76").Value = 500
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 1500
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("M76").Value = -1117

# CUL row 79
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H79").Value = 500
$ws.Range("I79").Value = 500
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 1500
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("M79").Value = -174

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 952.7143
$ws.Range("J131").Value = 1008.43665
$ws.Range("L131").Value = 3025.30995
$ws.Range("N131").Value = -13105.30995

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1166
$ws.Range("J135").Value = 1556.3334
$ws.Range("L135").Value = 14007.0006
$ws.Range("N135").Value = -19077.0006

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 9062.923000000001
$ws.Range("I137").Value = 11746.556
$ws.Range("J137").Value = 3024.75
$ws.Range("K137").Value = 35239.66800000001
$ws.Range("L137").Value = 9074.25
$ws.Range("M137").Value = -30139.66800000001
$ws.Range("N137").Value = -19274.25

# CUL row 138
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3867.8157
$ws.Range("I138").Value = 1203.3334
$ws.Range("J138").Value = 4367.4062
$ws.Range("K138").Value = 3610.0002
$ws.Range("L138").Value = 13102.2186
$ws.Range("M138").Value = 1529.9998
$ws.Range("N138").Value = -23382.2186

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1867.36
$ws.Range("J140").Value = 3671.4285
$ws.Range("L140").Value = 11014.2855
$ws.Range("N140").Value = -21374.2855

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4878.793
$ws.Range("I122").Value = 1940.7142
$ws.Range("J122").Value = 5813.636
$ws.Range("K122").Value = 5822.142599999999
$ws.Range("L122").Value = 17440.908
$ws.Range("M122").Value = -3372.142599999999
$ws.Range("N122").Value = -22340.908

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6749.5
$ws.Range("I132").Value = 10000
$ws.Range("J132").Value = 5666
$ws.Range("K132").Value = 30000
$ws.Range("L132").Value = 16998
$ws.Range("M132").Value = -27470
$ws.Range("N132").Value = -22058

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1158.25
$ws.Range("I46").Value = 1316.5
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 1316.5
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -1128.5
$ws.Range("N46").Value = -1376

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2997.5
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 2990
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 2990
$ws.Range("M100").Value = -2459
$ws.Range("N100").Value = -4072

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2891.9565
$ws.Range("I122").Value = 780
$ws.Range("J122").Value = 3478.611
$ws.Range("K122").Value = 2340
$ws.Range("L122").Value = 10435.833
$ws.Range("M122").Value = 110
$ws.Range("N122").Value = -15335.833

# WVR row 138
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 71097.60000000001
$ws.Range("J138").Value = 71097.60000000001
$ws.Range("L138").Value = 71097.60000000001
$ws.Range("N138").Value = -81377.60000000001
